# Natmi following Dr Hou advice
# Rewrites the Clec11a-Itga11 LR-pair table (rows 2-7) with updated
# per-cluster-pair statistics, adding the "ECs" target cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Clec11a"
$ws.Cells.Item(2,3).Value = "Itga11"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 9.930652333333333
$ws.Cells.Item(2,8).Value = 29.791957
$ws.Cells.Item(2,9).Value = 0.9673539331442913
$ws.Cells.Item(2,10).Value = 0.9673539331442912
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.04023833333333333
$ws.Cells.Item(2,14).Value = 0.120715
$ws.Cells.Item(2,15).Value = 0.001153919673903629
$ws.Cells.Item(2,16).Value = 0.001153919673903629
$ws.Cells.Item(2,17).Value = 0.3995928988061111
$ws.Cells.Item(2,18).Value = 3.596336089255
$ws.Cells.Item(2,19).Value = 0.001116248735083253
$ws.Cells.Item(2,20).Value = 0.001116248735083253

# Row 3
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Clec11a"
$ws.Cells.Item(3,3).Value = "Itga11"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 9.930652333333333
$ws.Cells.Item(3,8).Value = 29.791957
$ws.Cells.Item(3,9).Value = 0.9673539331442913
$ws.Cells.Item(3,10).Value = 0.9673539331442912
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 33.94639966666666
$ws.Cells.Item(3,14).Value = 101.839199
$ws.Cells.Item(3,15).Value = 0.9734851120464462
$ws.Cells.Item(3,16).Value = 0.9734851120464463
$ws.Cells.Item(3,17).Value = 337.1098930580492
$ws.Cells.Item(3,18).Value = 3033.989037522443
$ws.Cells.Item(3,19).Value = 0.9417046519955408
$ws.Cells.Item(3,20).Value = 0.9417046519955408

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Clec11a"
$ws.Cells.Item(4,3).Value = "Itga11"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 9.930652333333333
$ws.Cells.Item(4,8).Value = 29.791957
$ws.Cells.Item(4,9).Value = 0.9673539331442913
$ws.Cells.Item(4,10).Value = 0.9673539331442912
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.8843623333333334
$ws.Cells.Item(4,14).Value = 2.653087
$ws.Cells.Item(4,15).Value = 0.02536096827965006
$ws.Cells.Item(4,16).Value = 0.02536096827965006
$ws.Cells.Item(4,17).Value = 8.782294869028778
$ws.Cells.Item(4,18).Value = 79.040653821259
$ws.Cells.Item(4,19).Value = 0.0245330324136671
$ws.Cells.Item(4,20).Value = 0.02453303241366709

# Row 5
$ws.Cells.Item(5,1).Value = "sCs"
$ws.Cells.Item(5,2).Value = "Clec11a"
$ws.Cells.Item(5,3).Value = "Itga11"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.3351376666666666
$ws.Cells.Item(5,8).Value = 1.005413
$ws.Cells.Item(5,9).Value = 0.03264606685570878
$ws.Cells.Item(5,10).Value = 0.03264606685570878
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.04023833333333333
$ws.Cells.Item(5,14).Value = 0.120715
$ws.Cells.Item(5,15).Value = 0.001153919673903629
$ws.Cells.Item(5,16).Value = 0.001153919673903629
$ws.Cells.Item(5,17).Value = 0.01348538114388889
$ws.Cells.Item(5,18).Value = 0.121368430295
$ws.Cells.Item(5,19).Value = (3.767093882037554 / 100000)
$ws.Cells.Item(5,20).Value = (3.767093882037554 / 100000)

# Row 6
$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "Clec11a"
$ws.Cells.Item(6,3).Value = "Itga11"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.3351376666666666
$ws.Cells.Item(6,8).Value = 1.005413
$ws.Cells.Item(6,9).Value = 0.03264606685570878
$ws.Cells.Item(6,10).Value = 0.03264606685570878
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 33.94639966666666
$ws.Cells.Item(6,14).Value = 101.839199
$ws.Cells.Item(6,15).Value = 0.9734851120464462
$ws.Cells.Item(6,16).Value = 0.9734851120464463
$ws.Cells.Item(6,17).Value = 11.37671717602078
$ws.Cells.Item(6,18).Value = 102.390454584187
$ws.Cells.Item(6,19).Value = 0.03178046005090544
$ws.Cells.Item(6,20).Value = 0.03178046005090544

# Row 7
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Clec11a"
$ws.Cells.Item(7,3).Value = "Itga11"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.3351376666666666
$ws.Cells.Item(7,8).Value = 1.005413
$ws.Cells.Item(7,9).Value = 0.03264606685570878
$ws.Cells.Item(7,10).Value = 0.03264606685570878
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.8843623333333334
$ws.Cells.Item(7,14).Value = 2.653087
$ws.Cells.Item(7,15).Value = 0.02536096827965006
$ws.Cells.Item(7,16).Value = 0.02536096827965006
$ws.Cells.Item(7,17).Value = 0.2963831288812222
$ws.Cells.Item(7,18).Value = 2.667448159931
$ws.Cells.Item(7,19).Value = 0.0008279358659829657
$ws.Cells.Item(7,20).Value = 0.0008279358659829655
